# Final updates before printing!
# - Adjust the coordinates recorded for the "San Juans" kayaking entry (row 85)
# - Add three new visited-place rows (Boston Harbor, Pumice Plain, Nisqually
#   Wildlife Refuge) at the bottom of the table
# - Update the sheet view (zoom / frozen-pane scroll position / selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 85: corrected x/y coordinates ------------------------------------
$ws.Range("A85").Value = -123.094637
$ws.Range("B85").Value = 48.606879

# --- New row 102: Boston Harbor --------------------------------------------
$ws.Range("A102").Value = -122.903856
$ws.Range("B102").Value = 47.15527
$ws.Range("C102").Value = "Boston Harbor"
$ws.Range("D102").Value = "kayaking"
$ws.Range("E102").Value = "freshman"
$ws.Range("F102").Value = 42309
$ws.Range("F102").NumberFormat = "[$-409]d/mmm/yy;@"
$ws.Range("H102").Value = "kate"

# --- New row 103: Pumice Plain ----------------------------------------------
$ws.Range("A103").Value = -122.173267
$ws.Range("B103").Value = 46.242195
$ws.Range("C103").Value = "Pumice Plain"
$ws.Range("D103").Value = "hiking"
$ws.Range("E103").Value = "sophomore"
$ws.Range("F103").Value = 42583
$ws.Range("F103").NumberFormat = "[$-409]d/mmm/yy;@"
$ws.Range("H103").Value = "kate"

# --- New row 104: Nisqually Wildlife Refuge ---------------------------------
$ws.Range("A104").Value = -122.694864
$ws.Range("B104").Value = 47.082683
$ws.Range("C104").Value = "Nisqually Wildlife Refuge"
$ws.Range("D104").Value = "hiking"
$ws.Range("E104").Value = "senior"
$ws.Range("F104").Value = 43466
$ws.Range("F104").NumberFormat = "[$-409]d/mmm/yy;@"
$ws.Range("H104").Value = "kate"

# --- View: zoom back to 100%, keep header frozen, move selection -----------
$excel.ActiveWindow.Zoom = 100
$excel.ActiveWindow.FreezePanes = $false
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
[void]$ws.Range("F105").Select()
